$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number for every data row
# (rows 2 through 394). The commit bumps that date by one day (46081 -> 46082)
# for every row in the sheet.
$ws.Range("C2:C394").Value = 46082
